$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (233) down into
# the new rows (234:238) so the new date cells in column A keep the same
# date number format / style as the rest of the column.
$ws.Range("A233").Copy() | Out-Null
$ws.Range("A234:A238").PasteSpecial(-4122) | Out-Null

# New data rows
$ws.Range("A234").Value = 44308
$ws.Range("B234").Value = 8
$ws.Range("C234").Value = 43
$ws.Range("D234").Value = 239.1812214929358

$ws.Range("A235").Value = 44309
$ws.Range("B235").Value = 8
$ws.Range("C235").Value = 43
$ws.Range("D235").Value = 239.1812214929358

$ws.Range("A236").Value = 44310
$ws.Range("B236").Value = 9
$ws.Range("C236").Value = 47
$ws.Range("D236").Value = 261.430637445767

$ws.Range("A237").Value = 44311
$ws.Range("B237").Value = 10
$ws.Range("C237").Value = 50
$ws.Range("D237").Value = 278.1176994103905

$ws.Range("A238").Value = 44312
$ws.Range("B238").Value = 3
$ws.Range("C238").Value = 49
$ws.Range("D238").Value = 272.5553454221827
